# Read input data in txt form: append the two new "初始类别/ebsilon vs 聚类个数"
# result tables (epsilon-search results for two different kxi/cluster-count
# configurations) below the existing alpha-search table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Small spacer row (row 23) separating the existing table from the new ones.
$ws.Rows.Item(23).RowHeight = 5.25

# ---- First new block (rows 24-32) ----
$ws.Range("A24").Value = "初始类别"
$ws.Range("C24").Value = 40

$ws.Range("A25").Value = "ebsilon"
$ws.Range("C25").Value = "聚类个数"

$ws.Range("A26").Value = 0.001
$ws.Range("C26").Value = 25
$ws.Range("D26").Value = 29
$ws.Range("E26").Value = 34

$ws.Range("A27").Value = 0.01
$ws.Range("C27").Value = 31
$ws.Range("D27").Value = 29

$ws.Range("A28").Value = 0.1

$ws.Range("A29").Value = 1

$ws.Range("A30").Value = 10

$ws.Range("A31").Value = 100

$ws.Range("A32").Value = 1000
$ws.Range("C32").Value = 38
$ws.Range("D32").Value = 39
$ws.Range("E32").Value = 38

# ---- Second new block (rows 39-47) ----
$ws.Range("A39").Value = "初始类别"
$ws.Range("C39").Value = 40

$ws.Range("A40").Value = "ebsilon"
$ws.Range("C40").Value = "聚类个数"

$ws.Range("A41").Value = 0.001
$ws.Range("C41").Value = 18
$ws.Range("D41").Value = 19
$ws.Range("E41").Value = 18

$ws.Range("A42").Value = 0.01
$ws.Range("C42").Value = 19
$ws.Range("D42").Value = 19
$ws.Range("E42").Value = 19

$ws.Range("A43").Value = 0.1
$ws.Range("C43").Value = 18
$ws.Range("D43").Value = 20

$ws.Range("A44").Value = 1
$ws.Range("C44").Value = 19
$ws.Range("D44").Value = 19
$ws.Range("E44").Value = 20

$ws.Range("A45").Value = 10
$ws.Range("C45").Value = 19
$ws.Range("D45").Value = 19
$ws.Range("E45").Value = 19

$ws.Range("A46").Value = 100
$ws.Range("C46").Value = 21
$ws.Range("D46").Value = 18

$ws.Range("A47").Value = 1000

# Update the view: scroll so row 25 is at the top and select E46, matching
# where the user ended up after typing in the new data.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E46").Select()
